# Added comments and a second controller to the test hierarchy in the main()
#
# - Drop the trailing three "elapsed-time" columns (N, O, P) — the sample
#   table now only spans 12 data columns (B..M).
# - Update the existing CMD/RSP rows (3 & 4) so their state sequences read
#   INIT/INIT/INIT -> RUN x6 -> TERMINATE x3 (row 3) and
#   NOP -> EXECUTING (with two DONE blips) (row 4) over the shortened range.
# - Append a second controller pair to the hierarchy: a new
#   "CMD_parent-controller_to_child_controller" row and a new
#   "RSP_child_controller_to_parent-controller" row, each following the
#   same command/response state-machine shape as the first pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop columns N:P entirely (used range shrinks from A1:P4 to A1:M*).
$ws.Range("N1:P4").EntireColumn.Delete()

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M")

# Row 3: CMD_board-monitor_to_parent-controller
$row3vals = @("INIT","INIT","INIT","RUN","RUN","RUN","RUN","RUN","RUN","TERMINATE","TERMINATE","TERMINATE")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $row3vals[$i]
}

# Row 4: RSP_parent-controller_to_board-monitor
$row4vals = @("NOP","EXECUTING","EXECUTING","DONE","EXECUTING","EXECUTING","EXECUTING","EXECUTING","EXECUTING","DONE","EXECUTING","EXECUTING")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $row4vals[$i]
}

# Row 5 (new): CMD_parent-controller_to_child_controller
# Copy formatting (bold + border + centered header style) from A4 so the
# new row-label cell matches the look of the existing labels.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "CMD_parent-controller_to_child_controller"
$row5vals = @("NOP","INIT","INIT","INIT","RUN","RUN","RUN","RUN","RUN","RUN","TERMINATE","TERMINATE")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $row5vals[$i]
}

# Row 6 (new): RSP_child_controller_to_parent-controller
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "RSP_child_controller_to_parent-controller"
$row6vals = @("NOP","EXECUTING","DONE","DONE","EXECUTING","EXECUTING","EXECUTING","EXECUTING","DONE","DONE","EXECUTING","DONE")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "6").Value = $row6vals[$i]
}
